# modif calcul ball launcher
# The coil wire diameter ("Diamètre trou bobine", row 6, cell E6) is updated
# from 8 mm to 10.5 mm. All the dependent formulas elsewhere in the sheet
# (rows 13, 14, 16, 18-22, 26-29, columns E and H:U) are plain formulas that
# will be recomputed automatically by the recalculation engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("E6").Value = 10.5

# Keep the active selection in sync with the edited workbook, as observed
# after the edit (selection moved to M36 on Feuil1).
$ws.Range("M36").Select()
